# Adds a "Celltype"/"Color" lookup table in columns F:G of sheet1.
# Column F already held a list of cell-type names in F1:F30; we shift
# that list down one row (to make room for a header row), insert four
# additional cell-type rows, and add a parallel "Color" column (G) with
# a hex color code for every cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Celltype", "Color"),
    @("CD4 naive", "#7FC97F"),
    @("CD4 TCM", "#BEAED4"),
    @("CD4 TEM", "#FDC086"),
    @("CD4 CTL", "#FFFF99"),
    @("CD4 Treg", "#386CB0"),
    @("CD4 proliferating", "#F0027F"),
    @("CD8 naive", "#BF5B17"),
    @("CD8 TCM", "#666666"),
    @("CD8 TEM", "#1B9E77"),
    @("CD8 CTL", "#D95F02"),
    @("CD8 proliferating", "#7570B3"),
    @("MAIT", "#E7298A"),
    @("NKT", "#66A61E"),
    @("DNT", "#E6AB02"),
    @("GDT", "#A6761D"),
    @("NK CD16-", "#666666"),
    @("NK CD16+", "#A6CEE3"),
    @("NK proliferating", "#1F78B4"),
    @("B naive", "#B2DF8A"),
    @("B intermediate", "#33A02C"),
    @("B memory", "#FB9A99"),
    @("Plasma", "#E31A1C"),
    @("Classical monocyte", "#FDBF6F"),
    @("Intermediate monocyte", "#FF7F00"),
    @("Non-classical monocyte", "#CAB2D6"),
    @("CDC1", "#6A3D9A"),
    @("CDC2", "#FFFF99"),
    @("PDC", "#B15928"),
    @("ASDC", "#FBB4AE"),
    @("Platelet", "#B3CDE3"),
    @("HSPC", "#CCEBC5"),
    @("Lin-", "#DECBE4"),
    @("Multiplet", "#FED9A6"),
    @("Dead/debris", "#FFFFCC")
)

$row = 1
foreach ($pair in $data) {
    $ws.Cells.Item($row, 6).Value = $pair[0]
    $ws.Cells.Item($row, 7).Value = $pair[1]
    $row++
}

# Column F now holds "Intermediate monocyte" (22 chars) as its longest
# entry, same as the pre-existing bestFit column H, so widen it to match.
$ws.Columns.Item(6).ColumnWidth = 21.5

$ws.Range("L23").Select() | Out-Null
